$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 379, pushing the existing 44399/44400 rows
# (old 379-384) down to 384-389. Excel copies formatting from the row
# above on insert, which is what gives the shifted rows their correct
# look; the new blank rows (379-383) are filled in below with the new
# weekly price data (week of 44448).
$ws.Rows.Item(379).Resize(5).EntireRow.Insert()

# Row 379
$ws.Cells.Item(379, 1).Value = 3
$ws.Cells.Item(379, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(379, 3).Value = "Coquimbo"
$ws.Cells.Item(379, 4).Value = 44448
$ws.Cells.Item(379, 5).Value = 5
$ws.Cells.Item(379, 6).Value = "Fruta"
$ws.Cells.Item(379, 7).Value = 100101
$ws.Cells.Item(379, 8).Value = "Berries"
$ws.Cells.Item(379, 9).Value = 100101007
$ws.Cells.Item(379, 10).Value = "Kiwi"
$ws.Cells.Item(379, 11).Value = "Hayward"
$ws.Cells.Item(379, 12).Value = "Especial"
$ws.Cells.Item(379, 13).Value = 64
$ws.Cells.Item(379, 14).Value = 10000
$ws.Cells.Item(379, 15).Value = 10000
$ws.Cells.Item(379, 16).Value = 10000
$ws.Cells.Item(379, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(379, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(379, 19).Value = 1000
$ws.Cells.Item(379, 20).Value = 10

# Row 380
$ws.Cells.Item(380, 1).Value = 3
$ws.Cells.Item(380, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(380, 3).Value = "Coquimbo"
$ws.Cells.Item(380, 4).Value = 44448
$ws.Cells.Item(380, 5).Value = 5
$ws.Cells.Item(380, 6).Value = "Fruta"
$ws.Cells.Item(380, 7).Value = 100101
$ws.Cells.Item(380, 8).Value = "Berries"
$ws.Cells.Item(380, 9).Value = 100101007
$ws.Cells.Item(380, 10).Value = "Kiwi"
$ws.Cells.Item(380, 11).Value = "Hayward"
$ws.Cells.Item(380, 12).Value = "Primera"
$ws.Cells.Item(380, 13).Value = 68
$ws.Cells.Item(380, 14).Value = 9000
$ws.Cells.Item(380, 15).Value = 9000
$ws.Cells.Item(380, 16).Value = 9000
$ws.Cells.Item(380, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(380, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(380, 19).Value = 900
$ws.Cells.Item(380, 20).Value = 10

# Row 381
$ws.Cells.Item(381, 1).Value = 3
$ws.Cells.Item(381, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(381, 3).Value = "Coquimbo"
$ws.Cells.Item(381, 4).Value = 44448
$ws.Cells.Item(381, 5).Value = 5
$ws.Cells.Item(381, 6).Value = "Fruta"
$ws.Cells.Item(381, 7).Value = 100101
$ws.Cells.Item(381, 8).Value = "Berries"
$ws.Cells.Item(381, 9).Value = 100101007
$ws.Cells.Item(381, 10).Value = "Kiwi"
$ws.Cells.Item(381, 11).Value = "Hayward"
$ws.Cells.Item(381, 12).Value = "Primera"
$ws.Cells.Item(381, 13).Value = 56
$ws.Cells.Item(381, 14).Value = 12000
$ws.Cells.Item(381, 15).Value = 12000
$ws.Cells.Item(381, 16).Value = 12000
$ws.Cells.Item(381, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(381, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(381, 19).Value = 800
$ws.Cells.Item(381, 20).Value = 15

# Row 382
$ws.Cells.Item(382, 1).Value = 3
$ws.Cells.Item(382, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(382, 3).Value = "Coquimbo"
$ws.Cells.Item(382, 4).Value = 44448
$ws.Cells.Item(382, 5).Value = 5
$ws.Cells.Item(382, 6).Value = "Fruta"
$ws.Cells.Item(382, 7).Value = 100101
$ws.Cells.Item(382, 8).Value = "Berries"
$ws.Cells.Item(382, 9).Value = 100101007
$ws.Cells.Item(382, 10).Value = "Kiwi"
$ws.Cells.Item(382, 11).Value = "Hayward"
$ws.Cells.Item(382, 12).Value = "Segunda"
$ws.Cells.Item(382, 13).Value = 60
$ws.Cells.Item(382, 14).Value = 8000
$ws.Cells.Item(382, 15).Value = 8000
$ws.Cells.Item(382, 16).Value = 8000
$ws.Cells.Item(382, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(382, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(382, 19).Value = 800
$ws.Cells.Item(382, 20).Value = 10

# Row 383
$ws.Cells.Item(383, 1).Value = 3
$ws.Cells.Item(383, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(383, 3).Value = "Coquimbo"
$ws.Cells.Item(383, 4).Value = 44448
$ws.Cells.Item(383, 5).Value = 5
$ws.Cells.Item(383, 6).Value = "Fruta"
$ws.Cells.Item(383, 7).Value = 100101
$ws.Cells.Item(383, 8).Value = "Berries"
$ws.Cells.Item(383, 9).Value = 100101007
$ws.Cells.Item(383, 10).Value = "Kiwi"
$ws.Cells.Item(383, 11).Value = "Hayward"
$ws.Cells.Item(383, 12).Value = "Segunda"
$ws.Cells.Item(383, 13).Value = 50
$ws.Cells.Item(383, 14).Value = 11000
$ws.Cells.Item(383, 15).Value = 11000
$ws.Cells.Item(383, 16).Value = 11000
$ws.Cells.Item(383, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(383, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(383, 19).Value = 733
$ws.Cells.Item(383, 20).Value = 15
